# Apply updated cronograma (schedule) data: reorder task blocks and add
# descriptive text for each task in column B, per the commit:
# "adicionando as ultimas anotacoes antes da apresentacao. grafico de
# cronograma tb foi ajustada."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r=2; A='Revisão'; B='Revisão de literatura, estudo, desenvolvimento de protocolos de análise'; C=42856; D=43069; E='#95BB72'; F='white' }
    @{ r=3; A='Revisão'; B='Revisão de literatura, estudo, desenvolvimento de protocolos de análise'; C=43101; D=43434; E='#95BB72'; F='white' }
    @{ r=4; A='Revisão'; B='Revisão de literatura, estudo, desenvolvimento de protocolos de análise'; C=43466; D=43799; E='#95BB72'; F='white' }
    @{ r=5; A='Revisão'; B='Revisão de literatura, estudo, desenvolvimento de protocolos de análise'; C=43831; D=44165; E='#95BB72'; F='white' }
    @{ r=6; A='Revisão'; B='Revisão de literatura, estudo, desenvolvimento de protocolos de análise'; C=44197; D=44499; E='#95BB72'; F='white' }
    @{ r=7; A='Revisão'; B='Revisão de literatura, estudo, desenvolvimento de protocolos de análise'; C=44562; D=44681; E='#95BB72'; F='white' }
    @{ r=8; A='Análise '; B=' Análise de dados'; C=43009; D=43069; E='#4B6043'; F='white' }
    @{ r=9; A='Análise '; B=' Análise de dados'; C=43101; D=43434; E='#4B6043'; F='white' }
    @{ r=10; A='Análise '; B=' Análise de dados'; C=43466; D=43799; E='#4B6043'; F='white' }
    @{ r=11; A='Análise '; B=' Análise de dados'; C=43831; D=44165; E='#4B6043'; F='white' }
    @{ r=12; A='Análise '; B=' Análise de dados'; C=44197; D=44499; E='#4B6043'; F='white' }
    @{ r=13; A='Análise '; B=' Análise de dados'; C=44562; D=44681; E='#4B6043'; F='white' }
    @{ r=14; A='Coleta '; B='Coleta  de dados dos ensaios'; C=42856; D=43040; E='#DDEAD1'; F='black' }
    @{ r=15; A='Coleta '; B='Coleta  de dados dos ensaios'; C=43101; D=43160; E='#DDEAD1'; F='black' }
    @{ r=16; A='Interpretação '; B='Interpretação dos resultados, escrita e submissão artigos'; C=43009; D=43069; E='#658354'; F='white' }
    @{ r=17; A='Interpretação '; B='Interpretação dos resultados, escrita e submissão artigos'; C=43101; D=43434; E='#658354'; F='white' }
    @{ r=18; A='Interpretação '; B='Interpretação dos resultados, escrita e submissão artigos'; C=43466; D=43799; E='#658354'; F='white' }
    @{ r=19; A='Interpretação '; B='Interpretação dos resultados, escrita e submissão artigos'; C=43831; D=44165; E='#658354'; F='white' }
    @{ r=20; A='Interpretação '; B='Interpretação dos resultados, escrita e submissão artigos'; C=44197; D=44499; E='#658354'; F='white' }
    @{ r=21; A='Interpretação '; B='Interpretação dos resultados, escrita e submissão artigos'; C=44562; D=44681; E='#658354'; F='white' }
    @{ r=22; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=43040; D=43069; E='#C7DDB5'; F='black' }
    @{ r=23; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=43252; D=43281; E='#C7DDB5'; F='black' }
    @{ r=24; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=43405; D=43434; E='#C7DDB5'; F='black' }
    @{ r=25; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=43617; D=43646; E='#C7DDB5'; F='black' }
    @{ r=26; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=43770; D=43799; E='#C7DDB5'; F='black' }
    @{ r=27; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=44136; D=44165; E='#C7DDB5'; F='black' }
    @{ r=28; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=43983; D=44012; E='#C7DDB5'; F='black' }
    @{ r=29; A='Resultados'; B='Apresentação dos resultados parciais obtidos em seminários e relatórios (anuais e final)'; C=44501; D=44530; E='#C7DDB5'; F='black' }
    @{ r=30; A='Sessão Pública '; B='Apresentação dos resultados obtidos sessão pública.'; C=44743; D=44772; E='#B3CF99'; F='black' }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.A
    $ws.Cells.Item($row.r, 2).Value = $row.B
    $ws.Cells.Item($row.r, 3).Value = $row.C
    $ws.Cells.Item($row.r, 4).Value = $row.D
    $ws.Cells.Item($row.r, 5).Value = $row.E
    $ws.Cells.Item($row.r, 6).Value = $row.F
}

# Restore the on-screen selection state saved with the workbook (the user
# had scrolled down to row 16 and left the cursor on B30 before saving).
$ws.Range("B30").Select()
